# Applies the LOM3257.xlsx restructuring described in the commit diff.
# Strategy:
#  1. Update row 10 (Objetivos) B/C from the long "Desenvolver..." text to the
#     first docente's name - these cells already exist with correct styles,
#     so a plain value assignment is enough.
#  2. Delete rows 12-28 entirely (this removes the old "Docentes responsaveis"
#     block, the "Programa resumido" / "Programa" / bibliography blocks, and
#     the trailing requirements rows) so we start from a clean slate below
#     row 11.
#  3. Re-create rows 12-23 with the new layout: labels shifted up, the
#     docente names interleaved with the label rows, and the two
#     "Requisito" lines now directly following "Requisitos:" without gaps.
#  4. Copy cell formatting from known-good template cells (so the new cells
#     get the right bold/wrap/color styles) and set the custom row heights
#     that the new layout requires.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: row 10 content swap (style is already correct) ---
$ws.Range("B10").Value2 = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value2 = "5840730 - Antonio Jefferson da Silva Machado"

# --- Step 2: drop everything from row 12 down ---
$ws.Range("A12:C28").EntireRow.Delete()

# --- Step 3: write the new values for rows 12-23 ---
$ws.Range("A12").Value2 = "Programa resumido:"
$ws.Range("B12").Value2 = "3682251 - Gabrielle Weber Martins"
$ws.Range("C12").Value2 = "3682251 - Gabrielle Weber Martins"

$ws.Range("A13").Value2 = "Short syllabus:"

$ws.Range("A14").Value2 = "Programa:"
$ws.Range("B14").Value2 = "7797767 - Viktor Pastoukhov"
$ws.Range("C14").Value2 = "7797767 - Viktor Pastoukhov"

$ws.Range("A15").Value2 = "Syllabus:"

$ws.Range("A16").Value2 = "Avaliação:"

$ws.Range("A17").Value2 = "Método:"
$ws.Range("B17").Value2 = "5729033 - Weiliang Qian"
$ws.Range("C17").Value2 = "5729033 - Weiliang Qian"

$ws.Range("A18").Value2 = "Critério:"
$ws.Range("B18").Value2 = "A avaliação será composta por duas provas escritas (P1 e P2)."
$ws.Range("C18").Value2 = "A avaliação será composta por duas provas escritas (P1 e P2)."

$ws.Range("A19").Value2 = "Norma de recuperação:"
$ws.Range("B19").Value2 = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."
$ws.Range("C19").Value2 = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."

$ws.Range("A20").Value2 = "Bibliografia:"
$ws.Range("B20").Value2 = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C20").Value2 = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

$ws.Range("A21").Value2 = "Requisitos:"

$ws.Range("B22").Value2 = "LOB1004 -  Cálculo II  (Requisito)`n"
$ws.Range("C22").Value2 = "LOB1004 -  Cálculo II  (Requisito)`n"

$ws.Range("B23").Value2 = "LOB1018 -  Física I  (Requisito)`n"
$ws.Range("C23").Value2 = "LOB1018 -  Física I  (Requisito)`n"

# --- Step 4a: fix up cell formatting on the freshly created rows by copying
#     formats from template cells that already have the correct style
#     (A column = bold label style, B/C columns = wrap-text / red-wrap) ---
$ws.Range("A3").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C12").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A15").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C19").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A21").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Step 4b: row heights for the new layout (every row from 12-23 was
#     freshly (re)created by the delete above, so each custom height has to
#     be (re)applied explicitly; rows 16/21 keep the default height) ---
$ws.Rows(12).RowHeight = 60
$ws.Rows(13).RowHeight = 60
$ws.Rows(14).RowHeight = 120
$ws.Rows(15).RowHeight = 120
$ws.Rows(17).RowHeight = 60
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 60
$ws.Rows(20).RowHeight = 120
$ws.Rows(22).RowHeight = 30
$ws.Rows(23).RowHeight = 30
